$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to remain a text value even when the string looks
    # numeric (e.g. "0.999", "7.20"), matching the source data which
    # stores these as plain text/inline strings, not numbers.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "70.971.05"
$ws.Range("E2").Value = "  +2.39%  "
Set-TextValue $ws.Range("D3") "3.822.82"
$ws.Range("E3").Value = "  +1.10%  "
$ws.Range("E4").Value = "  -0.01%  "
Set-TextValue $ws.Range("D5") "669.47"
$ws.Range("E5").Value = "  +7.29%  "
Set-TextValue $ws.Range("D6") "169.86"
$ws.Range("E6").Value = "  +2.64%  "
Set-TextValue $ws.Range("D7") "3.819.21"
$ws.Range("E7").Value = "  +1.05%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("E9").Value = "  +1.24%  "
Set-TextValue $ws.Range("D10") "0.162"
$ws.Range("E10").Value = "  +0.92%  "
$ws.Range("E11").Value = "  +2.86%  "
Set-TextValue $ws.Range("D12") "7.03"
$ws.Range("E12").Value = "  +5.29%  "
Set-TextValue $ws.Range("D13") "0.0000246"
$ws.Range("E13").Value = "  -0.29%  "
Set-TextValue $ws.Range("D14") "36.22"
$ws.Range("E14").Value = "  +1.69%  "
Set-TextValue $ws.Range("D15") "4.466.97"
$ws.Range("E15").Value = "  +1.15%  "
Set-TextValue $ws.Range("D16") "3.811.44"
$ws.Range("E16").Value = "  -0.82%  "
Set-TextValue $ws.Range("D17") "70.873.13"
$ws.Range("E17").Value = "  +2.29%  "
Set-TextValue $ws.Range("D18") "17.82"
$ws.Range("E18").Value = "  +0.96%  "
$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue $ws.Range("D19") "11.76"
$ws.Range("E19").Value = "  +22.03%  "
$ws.Range("B20").Value = "Polkadot"
$ws.Range("C20").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue $ws.Range("D20") "7.20"
$ws.Range("E20").Value = "  +1.26%  "
$ws.Range("B21").Value = "TRON"
$ws.Range("C21").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue $ws.Range("D21") "0.114"
$ws.Range("E21").Value = "  +0.76%  "
Set-TextValue $ws.Range("D22") "476.95"
$ws.Range("E22").Value = "  +1.86%  "
Set-TextValue $ws.Range("D23") "0.718"
$ws.Range("E23").Value = "  +2.07%  "
Set-TextValue $ws.Range("D24") "83.28"
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("E25").Value = "  -2.87%  "
Set-TextValue $ws.Range("D26") "12.28"
$ws.Range("E26").Value = "  +2.23%  "
Set-TextValue $ws.Range("D27") "10.40"
$ws.Range("E27").Value = "  +3.73%  "
Set-TextValue $ws.Range("D28") "2.14"
$ws.Range("E28").Value = "  -1.10%  "
Set-TextValue $ws.Range("D30") "3.973.58"
$ws.Range("E30").Value = "  +1.09%  "
Set-TextValue $ws.Range("D31") "2.87"
$ws.Range("E31").Value = "  +7.70%  "
Set-TextValue $ws.Range("D32") "2.32"
$ws.Range("E32").Value = "  +3.72%  "
Set-TextValue $ws.Range("D33") "7.46"
$ws.Range("E33").Value = "  +2.50%  "
Set-TextValue $ws.Range("D34") "29.92"
$ws.Range("E34").Value = "  +3.94%  "
Set-TextValue $ws.Range("D35") "0.177"
$ws.Range("E35").Value = "  +7.27%  "
Set-TextValue $ws.Range("D36") "9.23"
$ws.Range("E36").Value = "  +2.67%  "
Set-TextValue $ws.Range("D37") "3.774.70"
$ws.Range("E37").Value = "  +1.14%  "
$ws.Range("E38").Value = "  -0.13%  "
$ws.Range("E39").Value = "  +0.73%  "
Set-TextValue $ws.Range("D40") "3.49"
$ws.Range("E40").Value = "  +2.81%  "
Set-TextValue $ws.Range("D41") "6.03"
$ws.Range("E41").Value = "  +3.72%  "
$ws.Range("E42").Value = "  +0.03%  "
Set-TextValue $ws.Range("D43") "0.999"
$ws.Range("E43").Value = "  -0.04%  "
$ws.Range("E44").Value = "  +10.33%  "
Set-TextValue $ws.Range("D46") "45.60"
$ws.Range("E46").Value = "  +4.96%  "
Set-TextValue $ws.Range("D47") "157.82"
$ws.Range("E47").Value = "  +3.68%  "
Set-TextValue $ws.Range("D48") "48.14"
$ws.Range("E48").Value = "  +2.93%  "
$ws.Range("E49").Value = "  +1.00%  "
Set-TextValue $ws.Range("D50") "1.43"
$ws.Range("E50").Value = "  +4.48%  "
Set-TextValue $ws.Range("D51") "8.54"
$ws.Range("E51").Value = "  +1.28%  "
